# Update profit files after running on 2025-10-23
# Append a new row (52) with the latest date/BTC/KAS allocation values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to be written as literal text so the date string is not
# auto-converted into a date serial number by Excel's smart entry parsing.
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = "10/23/2025"
# Drop the temporary text number format again so the new cell ends up with
# the default (unstyled) formatting, matching the rest of the data rows.
$ws.Range("A52").ClearFormats()

$ws.Range("B52").Value = 0.196346158260441
$ws.Range("C52").Value = 0.803653841739559
